$wb = $excel.ActiveWorkbook

$wsAnalysisUnit = $wb.Worksheets.Item("Analysis_Unit")
$wsVariable = $wb.Worksheets.Item("r AnalysisUnit_Variable")

# -----------------------------------------------------------------
# Sheet "Analysis_Unit": add two new columns (J: TeradataSchema,
# K: TargetVariable) and update the existing formula cell I3.
# -----------------------------------------------------------------

# Row 1 (style-only cells, same red header formatting as the rest of row 1)
$wsAnalysisUnit.Range("I1").Copy()
$wsAnalysisUnit.Range("J1:K1").PasteSpecial(-4122)

# Row 2 (header labels, same red header formatting as the rest of row 2)
$wsAnalysisUnit.Range("I2").Copy()
$wsAnalysisUnit.Range("J2:K2").PasteSpecial(-4122)
$wsAnalysisUnit.Range("J2").Value = "TeradataSchema"
$wsAnalysisUnit.Range("K2").Value = "TargetVariable"

# Row 3 (data values, plain formatting)
$wsAnalysisUnit.Range("J3").Value = "TEWSA0D"
$wsAnalysisUnit.Range("K3").Value = "TARGET"

# Update the existing external-info formula to take a parameter
$wsAnalysisUnit.Range("I3").Value = "LIB_EWS_IT.ExternalInfo(Param1);"

# New column J width
$wsAnalysisUnit.Columns.Item(10).ColumnWidth = 22.14

# Columns C and D now share the same (bestFit) width
$wsAnalysisUnit.Columns.Item(4).ColumnWidth = $wsAnalysisUnit.Columns.Item(3).ColumnWidth

# -----------------------------------------------------------------
# Sheet "r AnalysisUnit_Variable": append a new EXPOSURE variable row
# (row 20), mirroring the existing SEGMENT row (row 19).
# -----------------------------------------------------------------

$wsVariable.Range("A19:C19").Copy()
$wsVariable.Range("A20").PasteSpecial(-4122)
$wsVariable.Range("E19:F19").Copy()
$wsVariable.Range("E20").PasteSpecial(-4122)

$wsVariable.Range("A20").Value = "CREATE/MODIFY"
$wsVariable.Range("B20").Value = "CUSTOMER_EXPOSURE"
$wsVariable.Range("C20").Value = "CUSTOMER_EXPOSURE"
$wsVariable.Range("E20").Value = "CUSTOMER"
$wsVariable.Range("F20").Value = "EXPOSURE"

$wsVariable.Rows.Item(20).RowHeight = 15.75

# -----------------------------------------------------------------
# Selection / active-tab bookkeeping: "Analysis_Unit" becomes the
# active sheet (selected cell I3); "r AnalysisUnit_Variable" keeps a
# stored selection of C25 but is no longer the active tab.
# -----------------------------------------------------------------

$wsVariable.Activate()
$wsVariable.Range("C25").Select()

$wsAnalysisUnit.Activate()
$wsAnalysisUnit.Range("I3").Select()
